$d = $word.ActiveDocument

# Replace a whole paragraph's content by supplying well-formed OOXML for
# it (<w:pPr> plus the real run(s), deliberately omitting any originally
# "empty" leading run - the engine re-creates that empty run on its own).
# Going through InsertXML (rather than Range.Text / Find.Execute's
# replacement) avoids this runtime's Find/Replace behaviour of silently
# "smart-quoting" straight apostrophes in the replacement text.
function Set-ParagraphXml($ParaIndex, $InnerXml) {
    $p = $d.Paragraphs($ParaIndex)
    $rng = $d.Range($p.Range.Start, $p.Range.End)
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
           $InnerXml +
           '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

# Title (Heading1)
Set-ParagraphXml 1 '<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Play Aldo''s Journey for Free</w:t></w:r></w:p>'

# "What we like" bullets
Set-ParagraphXml 42 '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r><w:t>Exciting gameplay with unique bonuses in each region</w:t></w:r></w:p>'
Set-ParagraphXml 43 '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r><w:t>Four different regions to explore, offering diverse cultures and landscapes</w:t></w:r></w:p>'
Set-ParagraphXml 44 '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r><w:t>Well-defined and detailed graphics with simple symbols</w:t></w:r></w:p>'
Set-ParagraphXml 45 '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r><w:t>Moderate difficulty level, perfect for experienced gamblers</w:t></w:r></w:p>'

# "What we don't like" bullets
Set-ParagraphXml 47 '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r><w:t>Limited number of regions to explore</w:t></w:r></w:p>'
Set-ParagraphXml 48 '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r><w:t>Cartoonish graphics may not appeal to all players</w:t></w:r></w:p>'

# Closing bold title repeat
Set-ParagraphXml 49 '<w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Play Aldo''s Journey for Free</w:t></w:r></w:p>'

# Closing italic summary
Set-ParagraphXml 50 '<w:p><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Aldo''s Journey and play for free with unique bonuses.</w:t></w:r></w:p>'
